# Move the per-row comment text that was entered in column M ("comms")
# into column S ("comms_internal") for every data row (2-131) on the
# "Slovenian" sheet. This mirrors the author cutting the M2:M131 column
# and pasting it into S2, so each cell's value (and its formatting)
# travels with it while the now-empty M cells keep their original look.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reflect the live selection the author had right before doing the
# cut/paste: column M, rows 2 through 131 selected with M2 active.
$ws.Range("M2:M131").Select()

# Cut the whole block in one shot and drop it onto S2 - Excel expands
# the destination to the same shape as the source, so this lands the
# values in S2:S131 exactly, clears out M2:M131, and carries the M
# formatting (style index 5) onto every destination cell.
$ws.Range("M2:M131").Cut($ws.Range("S2"))
